# Update results of 20160404 - 001 -> 004
# Appends 10 new rows (12-21) of run results to the "logs" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Common string values reused across the new rows.
$preprocess1   = 'convert unicode to ascii, remove multiple spaces, trim "space" and ",", convert to lower'
$modelDetails1 = '2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 1000'
$preprocess2   = 'convert to lower, convert unicode to ascii, trim "space" and ",", remove multiple spaces'
$modelDetails2 = '2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 2000'

$features = '8 features: #ascii/(#ascii+#digit+#punctuation), #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, #max_digit_skip_0_1, first_character_ascii, first_character_digit'
$model = 'Neuron Network'
$templateFilter = '0 filters: '

# Data for the new rows: Time, RunningTime(s), Preprocess, Features, Model, Model_Details, Test_Accuracy, Val_Accuracy, Template Filter, (blank header col J)
$rows = @(
    @{ Row=12; Time='20160405_130344'; RunTime=548.509;  Preprocess=$preprocess1; ModelDetails=$modelDetails1; Test=0.978;              Val=0.99009900990099;  J=0.28125 },
    @{ Row=13; Time='20160405_131253'; RunTime=564.974;  Preprocess=$preprocess1; ModelDetails=$modelDetails1; Test=0.976666666666667;  Val=0.99009900990099;  J=0.270833333333333 },
    @{ Row=14; Time='20160405_132218'; RunTime=573.628;  Preprocess=$preprocess1; ModelDetails=$modelDetails1; Test=0.974666666666667;  Val=0.993399339933993;  J=0.268041237113402 },
    @{ Row=15; Time='20160405_133151'; RunTime=580.107;  Preprocess=$preprocess1; ModelDetails=$modelDetails1; Test=0.980666666666667;  Val=0.993399339933993;  J=0.175257731958763 },
    @{ Row=16; Time='20160405_134131'; RunTime=579.66;   Preprocess=$preprocess1; ModelDetails=$modelDetails1; Test=0.976;              Val=0.99009900990099;  J=0.302083333333333 },
    @{ Row=17; Time='20160405_145013'; RunTime=1103.617; Preprocess=$preprocess2; ModelDetails=$modelDetails2; Test=0.99;               Val=0.993399339933993;  J=0.175257731958763 },
    @{ Row=18; Time='20160405_150836'; RunTime=1107.887; Preprocess=$preprocess2; ModelDetails=$modelDetails2; Test=0.992;              Val=0.993399339933993;  J=0.195876288659794 },
    @{ Row=19; Time='20160405_152704'; RunTime=1150.024; Preprocess=$preprocess2; ModelDetails=$modelDetails2; Test=0.988666666666667;  Val=0.99009900990099;  J=0.322916666666667 },
    @{ Row=20; Time='20160405_154614'; RunTime=1208.437; Preprocess=$preprocess2; ModelDetails=$modelDetails2; Test=0.992666666666667;  Val=0.99009900990099;  J=0.291666666666667 },
    @{ Row=21; Time='20160405_160623'; RunTime=1249.055; Preprocess=$preprocess2; ModelDetails=$modelDetails2; Test=0.994;              Val=0.993399339933993;  J=0.144329896907216 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Cells.Item($n, 1).Value = $r.Time
    $ws.Cells.Item($n, 2).Value = $r.RunTime
    $ws.Cells.Item($n, 3).Value = $r.Preprocess
    $ws.Cells.Item($n, 4).Value = $features
    $ws.Cells.Item($n, 5).Value = $model
    $ws.Cells.Item($n, 6).Value = $r.ModelDetails
    $ws.Cells.Item($n, 7).Value = $r.Test
    $ws.Cells.Item($n, 8).Value = $r.Val
    $ws.Cells.Item($n, 9).Value = $templateFilter
    $ws.Cells.Item($n, 10).Value = $r.J
}
